$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.374.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.156.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.83%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.604"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.26"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.30%  "

$ws.Range("E8").Value = "  -0.06%  "

$ws.Range("E9").Value = "  -4.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.17%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.81%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0996"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.06%  "

$ws.Range("E14").Value = "  -3.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.475.30"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.141.96"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.776"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "41.216.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000102"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.77"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -12.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "226.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.65%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("E27").Value = "  -6.29%  "

$ws.Range("E28").Value = "  -9.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.78%  "

$ws.Range("E30").Value = "  -1.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.99%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "32.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0766"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.99%  "

$ws.Range("E35").Value = "  -9.37%  "

$ws.Range("E36").Value = "  -3.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.103"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.83%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0296"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.28%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -10.78%  "

$ws.Range("E41").Value = "  -3.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.187"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.80%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.79%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0951"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.87%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "95.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.70%  "

$ws.Range("E48").Value = "  -3.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("E51").Value = "  -7.75%  "
